# Bugfix: four municipalities (Bathurst, Campbellton, Dieppe, Edmundston) were
# incorrectly excluded when parsing the raw 2020 data. Insert them back as new
# rows at the top of the 2020 block (rows 2018-2021), which pushes the
# remaining 2020 rows down by four rows (old 2018-2108 -> new 2022-2112).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows right before the first existing 2020 data row (old row 2018).
# This shifts all following rows (including the rest of the 2020 block) down by 4.
$ws.Rows("2018:2021").Insert()

# Row 2018: Bathurst
$ws.Range("A2018").Value = 2020
$ws.Range("B2018").Value = "Bathurst"
$ws.Range("C2018").Value = 19153436
$ws.Range("D2018").Value = 3624702
$ws.Range("E2018").Value = 1013496
$ws.Range("F2018").Value = 1164400
$ws.Range("G2018").Value = 472498
$ws.Range("H2018").Value = 0
$ws.Range("I2018").Value = 500000
$ws.Range("J2018").Value = 0
$ws.Range("K2018").Value = 25928532

# Row 2019: Campbellton
$ws.Range("A2019").Value = 2020
$ws.Range("B2019").Value = "Campbellton"
$ws.Range("C2019").Value = 10686193
$ws.Range("D2019").Value = 1889297
$ws.Range("E2019").Value = 156231
$ws.Range("F2019").Value = 1093136
$ws.Range("G2019").Value = 503816
$ws.Range("H2019").Value = 46500
$ws.Range("I2019").Value = 405200
$ws.Range("J2019").Value = 66867
$ws.Range("K2019").Value = 14847240

# Row 2020: Dieppe
$ws.Range("A2020").Value = 2020
$ws.Range("B2020").Value = "Dieppe"
$ws.Range("C2020").Value = 53323334
$ws.Range("D2020").Value = 1095459
$ws.Range("E2020").Value = 345000
$ws.Range("F2020").Value = 1697000
$ws.Range("G2020").Value = 978500
$ws.Range("H2020").Value = 0
$ws.Range("I2020").Value = 926441
$ws.Range("J2020").Value = 583007
$ws.Range("K2020").Value = 58948741

# Row 2021: Edmundston
$ws.Range("A2021").Value = 2020
$ws.Range("B2021").Value = "Edmundston"
$ws.Range("C2021").Value = 23720297
$ws.Range("D2021").Value = 5564820
$ws.Range("E2021").Value = 1718846
$ws.Range("F2021").Value = 1342789
$ws.Range("G2021").Value = 725500
$ws.Range("H2021").Value = 5000
$ws.Range("I2021").Value = 2420632
$ws.Range("J2021").Value = 12801
$ws.Range("K2021").Value = 35510685

# Resize the table / list object so it (and its autofilter) covers the four
# newly-added rows, matching the sheet's new dimension (A1:K2112).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K2112"))
